$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets("ALC")
# row 40
$ws.Range("H40").Value = 4500
$ws.Range("J40").Value = 4500
$ws.Range("L40").Value = 4500
$ws.Range("N40").Value = -4850
# row 115
$ws.Range("H115").Value = 11000
# row 135
$ws.Range("H135").Value = 1957.0769
$ws.Range("I135").Value = 2141.3
$ws.Range("K135").Value = 19271.7
$ws.Range("M135").Value = -16736.7
# row 137
$ws.Range("H137").Value = 2463.65
$ws.Range("I137").Value = 1641.826
$ws.Range("J137").Value = 3575.5293
$ws.Range("K137").Value = 4925.478
$ws.Range("L137").Value = 10726.5879
$ws.Range("M137").Value = -2375.478
$ws.Range("N137").Value = -15826.5879
# row 138
$ws.Range("H138").Value = 5180.2383
$ws.Range("J138").Value = 5180.2383
$ws.Range("L138").Value = 15540.7149
$ws.Range("N138").Value = -25820.7149

# ---- Sheet: ARM ----
$ws = $wb.Worksheets("ARM")
# row 26
$ws.Range("H26").Value = 1000
$ws.Range("I26").Value = 1000
$ws.Range("K26").Value = 1000
$ws.Range("M26").Value = -670
# row 27
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
# row 32
$ws.Range("H32").Value = 21666.166
$ws.Range("I32").Value = 21666.166
$ws.Range("K32").Value = 21666.166
$ws.Range("M32").Value = -21379.166
# row 36
$ws.Range("H36").Value = 1200
$ws.Range("I36").Value = 1200
$ws.Range("K36").Value = 1200
$ws.Range("M36").Value = -854
# row 74
$ws.Range("H74").Value = 2580.65
$ws.Range("I74").Value = 1706.3
$ws.Range("K74").Value = 1706.3
$ws.Range("M74").Value = -832.3
# row 77
$ws.Range("H77").Value = 2580.65
$ws.Range("I77").Value = 1706.3
$ws.Range("K77").Value = 8531.5
$ws.Range("M77").Value = -4163.5
# row 131
$ws.Range("H131").Value = 80450
$ws.Range("J131").Value = 80450
$ws.Range("L131").Value = 80450
$ws.Range("N131").Value = -90530

# ---- Sheet: BSM ----
$ws = $wb.Worksheets("BSM")
# row 86
$ws.Range("H86").Value = 3958.3333
$ws.Range("I86").Value = 3798.0833
$ws.Range("J86").Value = 4599.3335
$ws.Range("K86").Value = 3798.0833
$ws.Range("L86").Value = 4599.3335
$ws.Range("M86").Value = -2675.0833
$ws.Range("N86").Value = -6845.3335
# row 89
$ws.Range("H89").Value = 3958.3333
$ws.Range("I89").Value = 3798.0833
$ws.Range("J89").Value = 4599.3335
$ws.Range("K89").Value = 18990.4165
$ws.Range("L89").Value = 22996.6675
$ws.Range("M89").Value = -13374.4165
$ws.Range("N89").Value = -34228.6675
# row 134
$ws.Range("H134").Value = 5142.4287
$ws.Range("I134").Value = 4500
$ws.Range("K134").Value = 13500
$ws.Range("M134").Value = -10965
# row 138
$ws.Range("H138").Value = 106214.75
$ws.Range("J138").Value = 106214.75
$ws.Range("L138").Value = 106214.75
$ws.Range("N138").Value = -116494.75

# ---- Sheet: CRP ----
$ws = $wb.Worksheets("CRP")
# row 31
$ws.Range("H31").Value = 2468.4583
$ws.Range("I31").Value = 2468.4583
$ws.Range("J31").Value = 0
$ws.Range("K31").Value = 2468.4583
$ws.Range("L31").Value = 0
$ws.Range("M31").Value = -2173.4583
$ws.Range("N31").Value = ""
# row 34
$ws.Range("H34").Value = 2468.4583
$ws.Range("I34").Value = 2468.4583
$ws.Range("J34").Value = 0
$ws.Range("K34").Value = 2468.4583
$ws.Range("L34").Value = 0
$ws.Range("M34").Value = -2266.4583
$ws.Range("N34").Value = ""
# row 36
$ws.Range("H36").Value = 2500
$ws.Range("I36").Value = 2500
$ws.Range("K36").Value = 2500
$ws.Range("M36").Value = -2112
# row 40
$ws.Range("H40").Value = 2500
$ws.Range("I40").Value = 2500
$ws.Range("K40").Value = 2500
$ws.Range("M40").Value = -2340
# row 105
$ws.Range("H105").Value = 3412.889
$ws.Range("I105").Value = 2959.4285
$ws.Range("J105").Value = 5000
$ws.Range("K105").Value = 2959.4285
$ws.Range("L105").Value = 5000
$ws.Range("M105").Value = -1212.4285
$ws.Range("N105").Value = -8494
# row 122
$ws.Range("H122").Value = 2322.25
$ws.Range("I122").Value = 2157.8333
$ws.Range("J122").Value = 2815.5
$ws.Range("K122").Value = 6473.499899999999
$ws.Range("L122").Value = 8446.5
$ws.Range("M122").Value = -4023.499899999999
$ws.Range("N122").Value = -13346.5
# row 134
$ws.Range("H134").Value = 2718.125
$ws.Range("I134").Value = 2457.6667
$ws.Range("K134").Value = 7373.000100000001
$ws.Range("M134").Value = -4838.000100000001
# row 141
$ws.Range("H141").Value = 66234.38
$ws.Range("J141").Value = 64546.1
$ws.Range("L141").Value = 64546.1
$ws.Range("N141").Value = -74906.10000000001

# ---- Sheet: CUL ----
$ws = $wb.Worksheets("CUL")
# row 5
$ws.Range("H5").Value = 729.125
$ws.Range("I5").Value = 738.8333
$ws.Range("K5").Value = 2216.4999
$ws.Range("M5").Value = -2104.4999
# row 44
$ws.Range("H44").Value = 33333508
$ws.Range("J44").Value = 442
$ws.Range("L44").Value = 1326
$ws.Range("N44").Value = -2122
# row 55
$ws.Range("H55").Value = 45454790
$ws.Range("I55").Value = 45454790
$ws.Range("J55").Value = 0
$ws.Range("K55").Value = 136364370
$ws.Range("L55").Value = 0
$ws.Range("M55").Value = -136364193
$ws.Range("N55").Value = ""
# row 107
$ws.Range("H107").Value = 316.33334
$ws.Range("J107").Value = 316.33334
$ws.Range("L107").Value = 949.0000200000001
$ws.Range("N107").Value = -4789.00002
# row 113
$ws.Range("H113").Value = 1997.25
$ws.Range("J113").Value = 1997.25
$ws.Range("L113").Value = 5991.75
$ws.Range("N113").Value = -10331.75
# row 131
$ws.Range("H131").Value = 2042.7222
$ws.Range("J131").Value = 2123.125
$ws.Range("L131").Value = 6369.375
$ws.Range("N131").Value = -16449.375
# row 135
$ws.Range("H135").Value = 729.125
$ws.Range("I135").Value = 738.8333
$ws.Range("K135").Value = 6649.4997
$ws.Range("M135").Value = -4114.4997
# row 141
$ws.Range("H141").Value = 6545.6
$ws.Range("I141").Value = 6545.6
$ws.Range("K141").Value = 19636.8
$ws.Range("M141").Value = -14456.8

# ---- Sheet: GSM ----
$ws = $wb.Worksheets("GSM")
# row 132
$ws.Range("H132").Value = 3282.4285
$ws.Range("J132").Value = 3282.4285
$ws.Range("L132").Value = 9847.2855
$ws.Range("N132").Value = -14907.2855

# ---- Sheet: LTW ----
$ws = $wb.Worksheets("LTW")
# row 7
$ws.Range("H7").Value = 1744
$ws.Range("J7").Value = 1743.5
$ws.Range("L7").Value = 1743.5
$ws.Range("N7").Value = -1967.5
# row 46
$ws.Range("H46").Value = 439.85715
$ws.Range("I46").Value = 439.85715
$ws.Range("K46").Value = 439.85715
$ws.Range("M46").Value = -251.85715
# row 93
$ws.Range("H93").Value = 2687.4
$ws.Range("I93").Value = 2741.6155
$ws.Range("K93").Value = 2741.6155
$ws.Range("M93").Value = -1493.6155
# row 100
$ws.Range("H100").Value = 1906.9286
$ws.Range("I100").Value = 1599.9
$ws.Range("J100").Value = 2674.5
$ws.Range("K100").Value = 1599.9
$ws.Range("L100").Value = 2674.5
$ws.Range("M100").Value = -1058.9
$ws.Range("N100").Value = -3756.5
# row 126
$ws.Range("H126").Value = 1744
$ws.Range("J126").Value = 1743.5
$ws.Range("L126").Value = 5230.5
$ws.Range("N126").Value = -10170.5

# ---- Sheet: WVR ----
$ws = $wb.Worksheets("WVR")
# row 122
$ws.Range("H122").Value = 2824.8462
$ws.Range("J122").Value = 3446.5
$ws.Range("L122").Value = 10339.5
$ws.Range("N122").Value = -15239.5
# row 132
$ws.Range("H132").Value = 4898.2
$ws.Range("I132").Value = 4500
$ws.Range("J132").Value = 4997.75
$ws.Range("K132").Value = 13500
$ws.Range("L132").Value = 14993.25
$ws.Range("M132").Value = -10970
$ws.Range("N132").Value = -20053.25
# row 137
$ws.Range("H137").Value = 87439.60000000001
$ws.Range("J137").Value = 87439.60000000001
$ws.Range("L137").Value = 87439.60000000001
$ws.Range("N137").Value = -97639.60000000001
